# Generate Report for Archive
# The localization status for this (single-row) report moved on from
# "Ready for handoff" to "In Translation" - update every sheet that
# surfaces that status column, then let Excel re-flow ("AutoFit") the
# columns whose content just got shorter so the sheet still looks right.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "zh-cn" (E) and "de-de" (F) status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Columns("E:F").AutoFit()

# --- Per-locale detail sheets: "Status" column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns("C:C").AutoFit()

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns("C:C").AutoFit()
